$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @('天奇股份', '三花智控', '华胜天成')
    3 = @('三花智控', '天奇股份', '三花智控')
    4 = @('华胜天成', '中国中铁', '天奇股份')
    5 = @('光线传媒', '贵州茅台', '风语筑')
    6 = @('利欧股份', '澜起科技', '博纳影业')
    7 = @('嘉美包装', '捷成股份', '利欧股份')
    8 = @('五洲新春', '百达精工', '克来机电')
    9 = @('万向钱潮', '东方财富', '光线传媒')
    10 = @('风语筑', '华胜天成', '嘉美包装')
    11 = @('百达精工', '嘉美包装', '协鑫集成')
    12 = @('博纳影业', '利亚德', '掌阅科技')
    13 = @('紫金矿业', '光线传媒', '首都在线')
    14 = @('中国中铁', '万向钱潮', '汉缆股份')
    15 = @('澜起科技', '利欧股份', '万向钱潮')
    16 = @('贵州茅台', '紫金矿业', '五洲新春')
    17 = @('掌阅科技', '卧龙电驱', '紫金矿业')
    18 = @('捷成股份', '蓝色光标', '大位科技')
    19 = @('汉缆股份', '兆易创新', '深科技')
    20 = @('蓝色光标', '风语筑', '二六三')
    21 = @('首都在线', '深科技', '蓝色光标')
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}